$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0. The trailing "_GoBack" bookmark currently sits at the end of the
#    "El tamano..." paragraph. Remove it now; we will re-create a
#    point bookmark at the end of the document once the new trailing
#    paragraph has been added.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 1. Insert a new, empty bold paragraph right before the
#    "El tamano del tipo de dato double..." paragraph.
# ------------------------------------------------------------------
$answerPara = $d.Paragraphs(2).Range
$answerPara.InsertParagraphBefore()
$blankPara = $d.Paragraphs(2).Range
$blankPara.Font.Bold = 1

# ------------------------------------------------------------------
# 2. Expand the wording inside the "El tamano..." paragraph
#    (now paragraph 3) to mention the Arduino Due / Arduino Uno sizes.
# ------------------------------------------------------------------
$rng = $d.Paragraphs(3).Range
$rng.Find.Execute(
    " es de 8-byte (64 bit),  el archivo ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " es de 8-byte (64 bit) en Arduino Due y en Arduino Uno ocupa 4 bytes,  el archivo ",
    2) | Out-Null

# Add the closing sentence about the float size, at the end of the
# paragraph (right after "tiene un tamano de ").
$rng2 = $d.Paragraphs(3).Range
$rng2.Collapse(0)
$rng2.InsertAfter("32 bits (4 bytes).")

# ------------------------------------------------------------------
# 3. Insert a brand-new paragraph after that one, explaining the
#    Arduino Uno case.
# ------------------------------------------------------------------
$endOfP3 = $d.Paragraphs(3).Range
$endOfP3.Collapse(0)
$endOfP3.InsertParagraphAfter()

$p4 = $d.Paragraphs(4).Range
$p4.Collapse(0)
# Append a trailing placeholder character so the text we are about to
# insert never sits exactly at the absolute end of the document while
# we create the bookmark (doing so at the literal end of the content
# mis-places the bookmark start marker).
$p4.InsertAfter("Para el Arduino Uno el cual vamos a manejar ambos tipos de datos tienen el mismo tamaño, con la diferencia que el float son decimales y los double son decimales con doble precisión.#")

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark as a zero-length point right at
#    the end of the document (i.e. right before the placeholder "#"),
#    then delete the placeholder character.
# ------------------------------------------------------------------
$bmPos = $p4.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($p4.End - 1, $p4.End)
$placeholder.Delete()
